$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-12-28"

# Update the label in A14
$ws.Range("A14").Value = "December (through 12-28)"

# December row (row 14) - updated arrest_made / no_arrest_made counts
$ws.Range("C14").Value = 38
$ws.Range("F14").Value = 83
$ws.Range("I14").Value = 96
$ws.Range("L14").Value = 63
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 52
$ws.Range("Q14").Value = 9
$ws.Range("R14").Value = 124
$ws.Range("U14").Value = 168

# Recompute arrest_rate columns for row 14 (arrest_made / (arrest_made + no_arrest_made))
$ws.Range("D14").Value = [Math]::Round($ws.Range("B14").Value2 / ($ws.Range("B14").Value2 + $ws.Range("C14").Value2), 4)
$ws.Range("G14").Value = [Math]::Round($ws.Range("E14").Value2 / ($ws.Range("E14").Value2 + $ws.Range("F14").Value2), 4)
$ws.Range("J14").Value = [Math]::Round($ws.Range("H14").Value2 / ($ws.Range("H14").Value2 + $ws.Range("I14").Value2), 4)
$ws.Range("M14").Value = [Math]::Round($ws.Range("K14").Value2 / ($ws.Range("K14").Value2 + $ws.Range("L14").Value2), 4)
$ws.Range("P14").Value = [Math]::Round($ws.Range("N14").Value2 / ($ws.Range("N14").Value2 + $ws.Range("O14").Value2), 4)
$ws.Range("S14").Value = [Math]::Round($ws.Range("Q14").Value2 / ($ws.Range("Q14").Value2 + $ws.Range("R14").Value2), 4)
$ws.Range("V14").Value = [Math]::Round($ws.Range("T14").Value2 / ($ws.Range("T14").Value2 + $ws.Range("U14").Value2), 4)

# Total row (row 15) - updated arrest_made / no_arrest_made counts
$ws.Range("C15").Value = 296
$ws.Range("F15").Value = 587
$ws.Range("I15").Value = 854
$ws.Range("L15").Value = 671
$ws.Range("N15").Value = 59
$ws.Range("O15").Value = 532
$ws.Range("Q15").Value = 73
$ws.Range("R15").Value = 1324
$ws.Range("U15").Value = 1710

# Recompute arrest_rate columns for row 15 (Total)
$ws.Range("D15").Value = [Math]::Round($ws.Range("B15").Value2 / ($ws.Range("B15").Value2 + $ws.Range("C15").Value2), 4)
$ws.Range("G15").Value = [Math]::Round($ws.Range("E15").Value2 / ($ws.Range("E15").Value2 + $ws.Range("F15").Value2), 4)
$ws.Range("J15").Value = [Math]::Round($ws.Range("H15").Value2 / ($ws.Range("H15").Value2 + $ws.Range("I15").Value2), 4)
$ws.Range("M15").Value = [Math]::Round($ws.Range("K15").Value2 / ($ws.Range("K15").Value2 + $ws.Range("L15").Value2), 4)
$ws.Range("P15").Value = [Math]::Round($ws.Range("N15").Value2 / ($ws.Range("N15").Value2 + $ws.Range("O15").Value2), 4)
$ws.Range("S15").Value = [Math]::Round($ws.Range("Q15").Value2 / ($ws.Range("Q15").Value2 + $ws.Range("R15").Value2), 4)
$ws.Range("V15").Value = [Math]::Round($ws.Range("T15").Value2 / ($ws.Range("T15").Value2 + $ws.Range("U15").Value2), 4)
